# QAPF with Result Report and Index Showing
#
# Applies to "Clastic.xlsx":
#  - Sheet "你的数据" (first tab): recolor the Marker "Color" column to "grey"
#    and drop the Alpha column to 0.4 for every existing data row, then add
#    three new "顶点" (vertex) rows (20-22) with the same marker formatting.
#  - Makes "你的数据" the active/selected tab (it was "Sheet2" before).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)   # "你的数据"

# --- Append three new "顶点" rows, copying the formatting of row 19 --------
# (added before the recolor pass below so the new shared strings land in the
# same table order as the saved workbook: "顶点" is registered before "grey")
$ws.Range("A17:K19").Copy()
$ws.Range("A20:K22").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("A20").Value = "顶点"
$ws.Range("B20").Value = 19
$ws.Range("C20").Value = 0
$ws.Range("D20").Value = 0
$ws.Range("E20").Value = 100
$ws.Range("F20").Value = "o"
$ws.Range("H20").Value = 50
$ws.Range("I20").Value = 0.4
$ws.Range("J20").Value = "-"
$ws.Range("K20").Value = 0.4

$ws.Range("A21").Value = "顶点"
$ws.Range("B21").Value = 20
$ws.Range("C21").Value = 0
$ws.Range("D21").Value = 100
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = "o"
$ws.Range("H21").Value = 50
$ws.Range("I21").Value = 0.4
$ws.Range("J21").Value = "-"
$ws.Range("K21").Value = 0.4

$ws.Range("A22").Value = "顶点"
$ws.Range("B22").Value = 21
$ws.Range("C22").Value = 100
$ws.Range("D22").Value = 0
$ws.Range("E22").Value = 0
$ws.Range("F22").Value = "o"
$ws.Range("H22").Value = 50
$ws.Range("I22").Value = 0.4
$ws.Range("J22").Value = "-"
$ws.Range("K22").Value = 0.4

# --- Recolor + re-alpha all 21 data rows (2..22) to match ------------------
$ws.Range("G2:G22").Value = "grey"
$ws.Range("I2:I22").Value = 0.4

# --- Make "你的数据" the active sheet / selection, as in the saved file ----
$ws.Activate()
$ws.Range("G16").Select()
